# Fruta / hortaliza, semanal
# The weekly refresh re-shuffles the date/volume/price/origin data for the
# "Mora" (blackberry) rows among themselves (rows 4 and 5 are untouched).
# Values below are the already-resolved target state for each row (D, M,
# N, O, P, R, S), taken straight from the new weekly extract.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($row, $d, $m, $n, $o, $p, $r, $s) {
    $ws.Cells.Item($row, 4).Value2  = $d   # D: Fecha
    $ws.Cells.Item($row, 13).Value2 = $m   # M: Volumen
    $ws.Cells.Item($row, 14).Value2 = $n   # N: Precio minimo
    $ws.Cells.Item($row, 15).Value2 = $o   # O: Precio maximo
    $ws.Cells.Item($row, 16).Value2 = $p   # P: Precio promedio ponderado
    $ws.Cells.Item($row, 18).Value2 = $r   # R: Origen
    $ws.Cells.Item($row, 19).Value2 = $s   # S: Precio $/Kg
}

Set-Row 2  44236 300 3600 4000 3800 "Provincia de Curicó" 1900
Set-Row 3  44237 100 3600 4000 3800 "Provincia de Curicó" 1900
Set-Row 6  44168 170 8000 8000 8000 "Provincia de Linares" 4000
Set-Row 7  44231 150 3400 3400 3400 "Provincia de Curicó" 1700
Set-Row 8  44208 85  3000 3000 3000 "Provincia de Linares" 1500
Set-Row 9  44232 200 3000 3000 3000 "Provincia de Curicó" 1500
Set-Row 10 44238 300 3600 4000 3800 "Provincia de Curicó" 1900
Set-Row 11 44188 150 3000 3400 3240 "Provincia de Linares" 1620
Set-Row 12 44194 120 3000 3000 3000 "Provincia de Linares" 1500
